# Updated symbol list on Mon Dec 19 11:06:43 UTC 2022 with GitHub Actions
# Refreshes price/volume/hour data for the coin listing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold numeric-looking text; force text
# number format first so values are written back as strings, matching
# the source inlineStr cells (preserves formats like "1.040", "0.04710").
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '248.26'
$ws.Range("G2").Value = '11'

$ws.Range("D3").Value = '21.65'
$ws.Range("G3").Value = '11'

$ws.Range("D4").Value = '5.501'
$ws.Range("G4").Value = '11'

$ws.Range("D5").Value = '0.05693'
$ws.Range("G5").Value = '11'

$ws.Range("D6").Value = '3.394'
$ws.Range("G6").Value = '11'

$ws.Range("D7").Value = '0.8041'
$ws.Range("G7").Value = '11'

$ws.Range("D8").Value = '1.040'
$ws.Range("G8").Value = '11'

$ws.Range("D9").Value = '0.1506'
$ws.Range("G9").Value = '11'

$ws.Range("D10").Value = '0.07619'
$ws.Range("G10").Value = '11'

$ws.Range("D11").Value = '0.03147'
$ws.Range("G11").Value = '11'

$ws.Range("G12").Value = '11'

$ws.Range("D13").Value = '0.09298'
$ws.Range("G13").Value = '11'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001653'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("G14").Value = '11'

$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '3.424'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").Value = '11'

$ws.Range("D16").Value = '0.04718'
$ws.Range("G16").Value = '11'

$ws.Range("D17").Value = '0.0005861'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("G17").Value = '11'

$ws.Range("D18").Value = '0.006357'
$ws.Range("G18").Value = '11'

$ws.Range("E19").Value = '18HotbitTokenHTBBestin24h'
$ws.Range("G19").Value = '11'

$ws.Range("D20").Value = '0.001043'
$ws.Range("G20").Value = '11'

$ws.Range("D21").Value = '0.0001501'
$ws.Range("G21").Value = '11'

$ws.Range("G22").Value = '11'

$ws.Range("D23").Value = '3.769'
$ws.Range("G23").Value = '11'

$ws.Range("D24").Value = '6.426'
$ws.Range("G24").Value = '11'

$ws.Range("D25").Value = '2.129'
$ws.Range("G25").Value = '11'

$ws.Range("D26").Value = '0.3280'
$ws.Range("G26").Value = '11'

$ws.Range("G27").Value = '11'

$ws.Range("G28").Value = '11'

$ws.Range("G29").Value = '11'

$ws.Range("G30").Value = '11'

$ws.Range("G31").Value = '11'

$ws.Range("G32").Value = '11'

$ws.Range("G33").Value = '11'

$ws.Range("G34").Value = '11'

$ws.Range("G35").Value = '11'

$ws.Range("G36").Value = '11'

$ws.Range("G37").Value = '11'

$ws.Range("G38").Value = '11'

$ws.Range("G39").Value = '11'

$ws.Range("D40").Value = '0.04114'
$ws.Range("G40").Value = '11'

$ws.Range("D41").Value = '0.006970'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("G41").Value = '11'

$ws.Range("D42").Value = '0.002972'
$ws.Range("G42").Value = '11'

$ws.Range("D43").Value = '0.1046'
$ws.Range("G43").Value = '11'

$ws.Range("D44").Value = '0.009138'
$ws.Range("G44").Value = '11'

$ws.Range("D45").Value = '0.00005839'
$ws.Range("G45").Value = '11'

$ws.Range("G46").Value = '11'

$ws.Range("D47").Value = '0.0005501'
$ws.Range("E47").Value = '46ACDXExchangeACXTWorstin24h'
$ws.Range("G47").Value = '11'

$ws.Range("D48").Value = '0.6825'
$ws.Range("G48").Value = '11'

$ws.Range("D49").Value = '0.008304'
$ws.Range("G49").Value = '11'

$ws.Range("G50").Value = '11'

$ws.Range("G51").Value = '11'
